# Update column G ("K") values on Sheet1 to the regenerated figures.
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
#  calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @(0, 1, 0, 2, 3, 0, 0, 0, 0, 1, 2, 1, 1, 1, 1, 0, 0, 2, 2, 0, 3, 3, 2, 1, 2, 1, 2, 0, 1, 2, 1, 1, 2, 2, 0, 2, 3, 1, 0, 1, 2, 1, 3, 2, 0, 1, 0, 1, 0, 0, 0, 0, 0, 2, 0, 1, 3, 3, 1, 1, 0, 2, 0, 0, 3, 1, 2, 2)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
